{"js": "// Update the date line and all 25 \"AxB=C\" multiplication answers in the\n// worksheet table with the new values from the target revision.\n// Every \"old\" text value in the document is unique, so a literal\n// search-and-replace (in document order) is sufficient and safe.\n\nconst replacements = [\n  [\"2024-03-06 Wednesday\", \"2024-03-07 Thursday\"],\n  [\"54\u00d732=1728\", \"32\u00d728=896\"],\n  [\"33\u00d777=2541\", \"62\u00d775=4650\"],\n  [\"59\u00d780=4720\", \"25\u00d736=900\"],\n  [\"48\u00d724=1152\", \"65\u00d786=5590\"],\n  [\"83\u00d764=5312\", \"62\u00d779=4898\"],\n  [\"70\u00d799=6930\", \"83\u00d731=2573\"],\n  [\"89\u00d744=3916\", \"69\u00d727=1863\"],\n  [\"27\u00d783=2241\", \"82\u00d774=6068\"],\n  [\"15\u00d771=1065\", \"45\u00d799=4455\"],\n  [\"16\u00d788=1408\", \"53\u00d759=3127\"],\n  [\"70\u00d772=5040\", \"45\u00d785=3825\"],\n  [\"81\u00d796=7776\", \"22\u00d798=2156\"],\n  [\"99\u00d778=7722\", \"84\u00d729=2436\"],\n  [\"81\u00d787=7047\", \"54\u00d792=4968\"],\n  [\"93\u00d755=5115\", \"96\u00d730=2880\"],\n  [\"79\u00d753=4187\", \"26\u00d733=858\"],\n  [\"97\u00d793=9021\", \"32\u00d757=1824\"],\n  [\"28\u00d720=560\", \"17\u00d738=646\"],\n  [\"89\u00d740=3560\", \"73\u00d760=4380\"],\n  [\"68\u00d778=5304\", \"18\u00d783=1494\"],\n  [\"91\u00d760=5460\", \"63\u00d796=6048\"],\n  [\"69\u00d739=2691\", \"17\u00d715=255\"],\n  [\"96\u00d765=6240\", \"42\u00d711=462\"],\n  [\"27\u00d775=2025\", \"34\u00d788=2992\"],\n  [\"77\u00d748=3696\", \"18\u00d783=1494\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const found = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  found.load(\"items\");\n  await context.sync();\n\n  if (found.items.length === 0) {\n    throw new Error(`Text not found: ${oldText}`);\n  }\n\n  found.items[0].insertText(newText, \"Replace\");\n}\n\nawait context.sync();\n", "ps1": "# Update the date line and all 25 \"AxB=C\" multiplication answers in the\n# worksheet table with the new values from the target revision.\n# Every \"old\" text value in the document is unique, so a plain\n# Find/Replace (one occurrence at a time, wdReplaceOne) is sufficient and\n# safe regardless of processing order.\n\n$d = $word.ActiveDocument\n\n$replacements = @(\n    @(\"2024-03-06 Wednesday\", \"2024-03-07 Thursday\"),\n    @(\"54\u00d732=1728\", \"32\u00d728=896\"),\n    @(\"33\u00d777=2541\", \"62\u00d775=4650\"),\n    @(\"59\u00d780=4720\", \"25\u00d736=900\"),\n    @(\"48\u00d724=1152\", \"65\u00d786=5590\"),\n    @(\"83\u00d764=5312\", \"62\u00d779=4898\"),\n    @(\"70\u00d799=6930\", \"83\u00d731=2573\"),\n    @(\"89\u00d744=3916\", \"69\u00d727=1863\"),\n    @(\"27\u00d783=2241\", \"82\u00d774=6068\"),\n    @(\"15\u00d771=1065\", \"45\u00d799=4455\"),\n    @(\"16\u00d788=1408\", \"53\u00d759=3127\"),\n    @(\"70\u00d772=5040\", \"45\u00d785=3825\"),\n    @(\"81\u00d796=7776\", \"22\u00d798=2156\"),\n    @(\"99\u00d778=7722\", \"84\u00d729=2436\"),\n    @(\"81\u00d787=7047\", \"54\u00d792=4968\"),\n    @(\"93\u00d755=5115\", \"96\u00d730=2880\"),\n    @(\"79\u00d753=4187\", \"26\u00d733=858\"),\n    @(\"97\u00d793=9021\", \"32\u00d757=1824\"),\n    @(\"28\u00d720=560\", \"17\u00d738=646\"),\n    @(\"89\u00d740=3560\", \"73\u00d760=4380\"),\n    @(\"68\u00d778=5304\", \"18\u00d783=1494\"),\n    @(\"91\u00d760=5460\", \"63\u00d796=6048\"),\n    @(\"69\u00d739=2691\", \"17\u00d715=255\"),\n    @(\"96\u00d765=6240\", \"42\u00d711=462\"),\n    @(\"27\u00d775=2025\", \"34\u00d788=2992\"),\n    @(\"77\u00d748=3696\", \"18\u00d783=1494\")\n)\n\nforeach ($pair in $replacements) {\n    $oldText = $pair[0]\n    $newText = $pair[1]\n\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    # wdFindContinue=1, Replace:=wdReplaceOne(1)\n    $find.Execute($oldText, $false, $false, $false, $false, $false, $true, 1, $false, $newText, 1) | Out-Null\n}\n\n$d.Save()\n"}
